# Auto-generated edit script: refresh market-derived Leve profit columns (H-N)
# per sheet, driven by updated market data (scheduled runner refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 120657
$ws.Range("J3").Value = 120657
$ws.Range("L3").Value = 120657
$ws.Range("N3").Value = -120885

$ws.Range("H9").Value = 6457.7
$ws.Range("I9").Value = 7511
$ws.Range("J9").Value = 4000
$ws.Range("K9").Value = 7511
$ws.Range("L9").Value = 4000
$ws.Range("M9").Value = -7342
$ws.Range("N9").Value = -4338

$ws.Range("H94").Value = 8000
$ws.Range("J94").Value = 8000
$ws.Range("L94").Value = 8000
$ws.Range("N94").Value = -8902

$ws.Range("H98").Value = 4481.758
$ws.Range("I98").Value = 4107
$ws.Range("K98").Value = 4107
$ws.Range("M98").Value = -2609

$ws.Range("H100").Value = 4873.375
$ws.Range("I100").Value = 2999
$ws.Range("K100").Value = 2999
$ws.Range("M100").Value = -2458

$ws.Range("H102").Value = 120657
$ws.Range("J102").Value = 120657
$ws.Range("L102").Value = 120657
$ws.Range("N102").Value = -127147

$ws.Range("H116").Value = 8267.733
$ws.Range("I116").Value = 8087.727
$ws.Range("J116").Value = 8371.947
$ws.Range("K116").Value = 8087.727
$ws.Range("L116").Value = 8371.947
$ws.Range("M116").Value = -4645.727
$ws.Range("N116").Value = -15255.947

$ws.Range("H122").Value = 4481.758
$ws.Range("I122").Value = 4107
$ws.Range("K122").Value = 12321
$ws.Range("M122").Value = -9871

$ws.Range("H129").Value = 2512.375
$ws.Range("I129").Value = 1715.6666
$ws.Range("J129").Value = 2990.4
$ws.Range("K129").Value = 5146.9998
$ws.Range("L129").Value = 8971.200000000001
$ws.Range("M129").Value = -146.9997999999996
$ws.Range("N129").Value = -18971.2

$ws.Range("H131").Value = 2317.1667
$ws.Range("I131").Value = 1846
$ws.Range("J131").Value = 7500
$ws.Range("K131").Value = 5538
$ws.Range("L131").Value = 22500
$ws.Range("M131").Value = -498
$ws.Range("N131").Value = -32580

$ws.Range("H138").Value = 6142.892
$ws.Range("I138").Value = 1627.1111
$ws.Range("J138").Value = 7872.3403
$ws.Range("K138").Value = 4881.3333
$ws.Range("L138").Value = 23617.0209
$ws.Range("M138").Value = 258.6666999999998
$ws.Range("N138").Value = -33897.0209


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2498.65
$ws.Range("I2").Value = 2041.0769
$ws.Range("J2").Value = 3348.4285
$ws.Range("K2").Value = 2041.0769
$ws.Range("L2").Value = 3348.4285
$ws.Range("M2").Value = -1928.0769
$ws.Range("N2").Value = -3574.4285

$ws.Range("H32").Value = 5520.452
$ws.Range("I32").Value = 5024.8057
$ws.Range("K32").Value = 5024.8057
$ws.Range("M32").Value = -4737.8057

$ws.Range("H61").Value = 45457504
$ws.Range("I61").Value = 45457504
$ws.Range("K61").Value = 45457504
$ws.Range("M61").Value = -45457292

$ws.Range("H102").Value = 3883.9167
$ws.Range("I102").Value = 3460.8
$ws.Range("K102").Value = 3460.8
$ws.Range("M102").Value = -1838.8

$ws.Range("H110").Value = 8453.429
$ws.Range("I110").Value = 8337.799999999999
$ws.Range("K110").Value = 8337.799999999999
$ws.Range("M110").Value = -6292.799999999999

$ws.Range("H116").Value = 2498.65
$ws.Range("I116").Value = 2041.0769
$ws.Range("J116").Value = 3348.4285
$ws.Range("K116").Value = 2041.0769
$ws.Range("L116").Value = 3348.4285
$ws.Range("M116").Value = 252.9231
$ws.Range("N116").Value = -7936.4285

$ws.Range("H132").Value = 29486020
$ws.Range("I132").Value = 4380.4346
$ws.Range("K132").Value = 13141.3038
$ws.Range("M132").Value = -10611.3038

$ws.Range("H136").Value = 45457504
$ws.Range("I136").Value = 45457504
$ws.Range("K136").Value = 136372512
$ws.Range("M136").Value = -136369962


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2498.65
$ws.Range("I3").Value = 2041.0769
$ws.Range("J3").Value = 3348.4285
$ws.Range("K3").Value = 2041.0769
$ws.Range("L3").Value = 3348.4285
$ws.Range("M3").Value = -1927.0769
$ws.Range("N3").Value = -3576.4285

$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").ClearContents()
$ws.Range("N53").Value = 0

$ws.Range("H86").Value = 8576.375
$ws.Range("I86").Value = 11851.083
$ws.Range("J86").Value = 5301.6665
$ws.Range("K86").Value = 11851.083
$ws.Range("L86").Value = 5301.6665
$ws.Range("M86").Value = -10728.083
$ws.Range("N86").Value = -7547.6665

$ws.Range("H89").Value = 8576.375
$ws.Range("I89").Value = 11851.083
$ws.Range("J89").Value = 5301.6665
$ws.Range("K89").Value = 59255.415
$ws.Range("L89").Value = 26508.3325
$ws.Range("M89").Value = -53639.415
$ws.Range("N89").Value = -37740.3325

$ws.Range("H94").Value = 1128.0286
$ws.Range("I94").Value = 907.2
$ws.Range("K94").Value = 907.2
$ws.Range("M94").Value = -456.2

$ws.Range("H99").Value = 4893.1
$ws.Range("I99").Value = 3474
$ws.Range("K99").Value = 3474
$ws.Range("M99").Value = -1976

$ws.Range("H107").Value = 5337.6665
$ws.Range("I107").Value = 15
$ws.Range("K107").Value = 15
$ws.Range("M107").Value = 1905


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5598.8
$ws.Range("I16").Value = 5999
$ws.Range("J16").Value = 4998.5
$ws.Range("K16").Value = 5999
$ws.Range("L16").Value = 4998.5
$ws.Range("M16").Value = -5712
$ws.Range("N16").Value = -5572.5

$ws.Range("H31").Value = 13163047
$ws.Range("I31").Value = 3205.1428
$ws.Range("J31").Value = 50010604
$ws.Range("K31").Value = 3205.1428
$ws.Range("L31").Value = 50010604
$ws.Range("M31").Value = -2910.1428
$ws.Range("N31").Value = -50011194

$ws.Range("H34").Value = 13163047
$ws.Range("I34").Value = 3205.1428
$ws.Range("J34").Value = 50010604
$ws.Range("K34").Value = 3205.1428
$ws.Range("L34").Value = 50010604
$ws.Range("M34").Value = -3003.1428
$ws.Range("N34").Value = -50011008

$ws.Range("H99").Value = 8004.273
$ws.Range("I99").Value = 7849.7
$ws.Range("J99").Value = 9550
$ws.Range("K99").Value = 7849.7
$ws.Range("L99").Value = 9550
$ws.Range("M99").Value = -6351.7
$ws.Range("N99").Value = -12546

$ws.Range("H113").Value = 5598.8
$ws.Range("I113").Value = 5999
$ws.Range("J113").Value = 4998.5
$ws.Range("K113").Value = 5999
$ws.Range("L113").Value = 4998.5
$ws.Range("M113").Value = -3829
$ws.Range("N113").Value = -9338.5

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").ClearContents()
$ws.Range("N118").Value = 0

$ws.Range("H126").Value = 8004.273
$ws.Range("I126").Value = 7849.7
$ws.Range("J126").Value = 9550
$ws.Range("K126").Value = 23549.1
$ws.Range("L126").Value = 28650
$ws.Range("M126").Value = -21079.1
$ws.Range("N126").Value = -33590


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 844
$ws.Range("J86").Value = 1491.3334
$ws.Range("L86").Value = 4474.0002
$ws.Range("N86").Value = -6846.0002

$ws.Range("H89").Value = 844
$ws.Range("J89").Value = 1491.3334
$ws.Range("L89").Value = 13422.0006
$ws.Range("N89").Value = -25278.0006

$ws.Range("H113").Value = 1659.6
$ws.Range("I113").Value = 300
$ws.Range("K113").Value = 900
$ws.Range("M113").Value = 1270

$ws.Range("H114").Value = 1000
$ws.Range("J114").Value = 1000
$ws.Range("L114").Value = 3000
$ws.Range("N114").Value = -9508


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 23000000
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H97").Value = 1250.7333
$ws.Range("I97").Value = 1553.7778
$ws.Range("K97").Value = 1553.7778
$ws.Range("M97").Value = -1057.7778

$ws.Range("H102").Value = 934.5833
$ws.Range("I102").Value = 837.8182
$ws.Range("K102").Value = 837.8182
$ws.Range("M102").Value = 784.1818

$ws.Range("H113").Value = 5666.3335
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 7499.5
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 7499.5
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -11839.5

$ws.Range("H126").Value = 54523.527
$ws.Range("I126").Value = 68243.8
$ws.Range("K126").Value = 204731.4
$ws.Range("M126").Value = -202261.4

$ws.Range("H132").Value = 8434.615
$ws.Range("I132").Value = 8605.454
$ws.Range("K132").Value = 25816.362
$ws.Range("M132").Value = -23286.362


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H68").Value = 2428.7334
$ws.Range("I68").Value = 2415.8
$ws.Range("K68").Value = 2415.8
$ws.Range("M68").Value = -1666.8

$ws.Range("H71").Value = 2428.7334
$ws.Range("I71").Value = 2415.8
$ws.Range("K71").Value = 12079
$ws.Range("M71").Value = -8335

$ws.Range("H100").Value = 3604
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5998.5
$ws.Range("I132").Value = 5998.5
$ws.Range("K132").Value = 17995.5
$ws.Range("M132").Value = -15465.5

$ws.Range("H136").Value = 3042.0952
$ws.Range("I136").Value = 2381.4707
$ws.Range("J136").Value = 5849.75
$ws.Range("K136").Value = 7144.4121
$ws.Range("L136").Value = 17549.25
$ws.Range("M136").Value = -4594.4121
$ws.Range("N136").Value = -22649.25

